$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 156.5814015774635
$ws.Range("C2").Value = 2977.411704809431
$ws.Range("D2").Value = 1671.891461163591
